# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# values to the latest scrape, as captured on
# Tue Aug  1 13:59:16 UTC 2023 by the GitHub Actions cron job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price column keeps being stored/read as plain text
# (matching the source data's dotted-thousands formatting), so Excel
# doesn't reinterpret values such as "245.00" or "1.750" as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.869.71"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.88"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.00"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07704"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3050"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.41"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.830.34"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.094"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6822"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.422"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008312"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.868.80"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.51"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.074.78"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  -3.70%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.94"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.24"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.220"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.153"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05093"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7754"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.142"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01849"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227.27"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9527"
$ws.Range("E41").Value = "  +6.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.05"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.900"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.606"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.975.64"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5157"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.14"
$ws.Range("E49").Value = "  -8.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.750"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.927"
$ws.Range("E51").Value = "  -1.14%  "
